# Updates the cryptos list on Sheet1 with refreshed price/volume figures,
# and fixes the relative ordering of a few coin rows (B/C columns) so that
# each row's name/link again matches its (re-ranked) price/volume data.
#
# Price values in column D are written as text (matching the source data,
# which uses dotted "thousands" groupings like "49.112.84" that are not
# valid numbers). Because some individual price strings do parse as plain
# decimals (e.g. "0.999", "111.72"), Excel would otherwise silently store
# them as numeric cells. To keep them as text - consistent with the rest
# of the column - NumberFormat is set to Text ("@") right before the
# assignment, and the cell style is reset back to "Normal" immediately
# after so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '49.112.84'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.17%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.616.51'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.02%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '111.72'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.50%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '322.70'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.93%  '

$ws.Range("E7").Value = '  -1.28%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.07%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.540'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.00%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.66'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.06%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.70'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.70%  '

$ws.Range("E12").Value = '  -1.11%  '

$ws.Range("E13").Value = '  +1.12%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.24'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.06%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.021.17'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.07%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.610.49'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.09%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.855'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.85%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '49.020.74'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.09%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.02'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.71%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.88'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.50%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.67'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.83%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0942'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.04%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '269.35'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.22%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '68.44'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.54%  '

$ws.Range("E25").Value = '  -1.66%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.12'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.22%  '

$ws.Range("E27").Value = '  +0.08%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.24'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.82%  '

$ws.Range("E29").Value = '  -0.46%  '

$ws.Range("B30").Value = 'Kaspa'
$ws.Range("C30").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.138'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.07%  '

$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '34.85'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.04%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.54'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.37%  '

$ws.Range("E33").Value = '  +0.28%  '

$ws.Range("E34").Value = '  +1.50%  '

$ws.Range("E35").Value = '  -0.12%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '19.00'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.72%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.93'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.04%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.03'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.65%  '

$ws.Range("E39").Value = '  +1.24%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '125.73'
$ws.Range("D40").Style = "Normal"

$ws.Range("E41").Value = '  -1.65%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '21.97'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.58%  '

$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0319'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.47%  '

$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.13'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.28%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.054.11'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.36%  '

$ws.Range("E46").Value = '  +7.54%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.20'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.50%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.12'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.75%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.90'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.70%  '

$ws.Range("B50").Value = 'MultiversX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '58.83'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.38%  '

$ws.Range("B51").Value = 'THORChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.18'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.99%  '
